$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the INSERT INTO ... shared text (used by A2:A11) to include the new legacy_column field
$newInsertText = "INSERT INTO datadictionary (``instance``,``scope``,``entity``,``prop_order``,``prop_name``,``prop_description``,``data_values``,``data_calculation``,``required``,``db_table``,``db_column``,``db_order``,``db_datatype``,``db_datalength``,``db_nullable``,``db_default``,``notes``,``legacy_column``) VALUES ('"
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("A$r").Value = $newInsertText
}

# Add new header for column U
$ws.Range("U1").Value = "legacy_column"

# Update the formulas in B2:B11 and C2:C11 to include the new U column reference
for ($r = 2; $r -le 11; $r++) {
    $ws.Range("B$r").Formula = "=D$r & `"','`" & E$r & `"','`" & F$r & `"','`" & G$r & `"','`" & H$r & `"','`" & I$r & `"','`" & J$r & `"','`" & K$r & `"','`" & L$r & `"','`" & M$r & `"','`" & N$r & `"','`" & O$r & `"','`" & P$r & `"','`" & Q$r & `"','`" & R$r & `"','`" & S$r & `"','`" & T$r & `"','`" & U$r & `"');`""
    $ws.Range("C$r").Formula = "=A$r & B$r"
}

$ws.Columns.Item(21).ColumnWidth = 12.5
